# Add new rows of gyroscope data captured on May 9th:
#  - 7 new rows are inserted at the top of the data block (old row 2 and
#    below shift down by 7 rows)
#  - 3 new rows are appended after the (now shifted) last existing row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 blank rows above the current data (old row 2..21 -> 9..28),
# then strip the formatting Excel auto-copies on insert so the new rows
# stay unstyled like the rest of the data rows.
$ws.Rows("2:8").Insert()
$ws.Rows("2:8").ClearFormats()

# Fill in the newly inserted rows with the earlier (May 9th) samples.
$newTop = @(
    @(-0.0421497002243995, 0.0813977941870689, -0.0167987942695617),
    @(-0.0009162978967650999, -0.0372627787292003, -0.051312681287527),
    @(-0.0261144898831844, 0.06383541971445079, -0.0560468845069408),
    @(0.0224492978304624, 0.0224492978304624, 0.0267253536731004),
    @(0.0284052342176437, -0.0120645882561802, 0.0143553335219621),
    @(0.0213802829384803, 0.0073303831741213, -0.0117591563612222),
    @(-0.0149661982432007, 0.042302418500185, 0.0258090570569038)
)

$r = 2
foreach ($row in $newTop) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $r = $r + 1
}

# Append 3 more (later) samples after the existing data, which is now at
# rows 9..28, so the new rows land at 29..31.
$newBottom = @(
    @(-0.050854530185461, -0.0387899428606033, -0.0474947728216648),
    @(-0.0007635815418325, 0.0126754539087414, 0.0435241498053073),
    @(0.027030786499381, 0.0493273697793483, -0.0366519130766391)
)

$r = 29
foreach ($row in $newBottom) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $r = $r + 1
}
